# weight_over_time.xlsx - add new (lowest) weight readings
# Commit: "add new lowest weight; y axis dual axis scale in decistones"
#
# The raw_data sheet gets 7 new rows (135-141) of date/time/weight
# observations, continuing the existing table (A:date, B:time, C:weight,
# D:=IF(B<TIME(12,0,0),"AM","PM")). We copy the formatting down from the
# last existing row (134) so the new cells pick up the same number
# formats (m/d/yyyy h:mm / h:mm) as the rest of the column, then fill in
# the values and the TOD formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# Extend formatting (number formats etc.) from the last populated row down
# across the new rows, matching how Excel extends a table when you type
# into the row below it.
$ws.Range("A134:D134").Copy()
$ws.Range("A135:D141").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New readings (most-recent-first, matching the existing sheet's order)
$dates   = @(44092.333333333336, 44092.333333333336, 44092.309027777781, 44091.888194444444, 44091.352083333331, 44091.351388888892, 44092.388194444444)
$times   = @(0.33333333333333331, 0.33333333333333331, 0.30902777777777779, 0.8881944444444444, 0.3520833333333333, 0.35138888888888892, 0.38819444444444445)
$weights = @(71.2, 71.900000000000006, 71.599999999999994, 72.3, 71.599999999999994, 71.599999999999994, 70.900000000000006)

for ($i = 0; $i -lt 7; $i++) {
    $r = 135 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $times[$i]
    $ws.Cells.Item($r, 3).Value = $weights[$i]
    $ws.Cells.Item($r, 4).Formula = "=IF(B" + $r + "<TIME(12,0,0), ""AM"", ""PM"")"
}

# Update the view to mirror where the user left off after entering the
# new rows.
$ws.Range("A141").Select()
$ws.Application.ActiveWindow.ScrollRow = 118
